# StudentSem.xlsx refactor: redistribute "Semester" numbering (column A) into
# even groups of 100 rows (was 150/150/150), fix the corrupted/duplicated
# StudentID block that used to restart at 24000022, and extend the sheet with
# 50 more StudentID rows (continuing the arithmetic +5 sequence) for the new
# 5th semester group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (Semester) renumbering -------------------------------------
# Rows 2-101   : Semester 1 (unchanged, already correct)
# Rows 102-151 : Semester 1 -> 2
$ws.Range("A102:A151").Value = 2
# Rows 152-201 : Semester 2 (unchanged, already correct)
# Rows 202-301 : Semester 2 -> 3
$ws.Range("A202:A301").Value = 3
# Rows 302-401 : Semester 3 -> 4
$ws.Range("A302:A401").Value = 4
# Rows 402-451 : Semester 3 -> 5
$ws.Range("A402:A451").Value = 5

# --- Column B (StudentID) fix + extension --------------------------------
# Rows 302 onward previously re-used stale/duplicate StudentIDs
# (restarting at 24000022). Replace them with the real continuation of the
# sequence that row 301 (24004507) was already following, stepping by 5,
# all the way through the 50 brand-new rows (452-501).
$b = 24004512
for ($r = 302; $r -le 501; $r++) {
    $ws.Cells.Item($r, 2).Value = $b
    $b = $b + 5
}

# --- New rows 452-501 need Semester + Programme filled in too ------------
$ws.Range("A452:A501").Value = 5
for ($r = 452; $r -le 501; $r++) {
    $ws.Cells.Item($r, 3).Value = "BCS"
}

# --- Selection / active cell, to mirror the author's final cursor spot ---
$null = $ws.Range("F193").Select()
